# Natmi following Dr Hou advice
#
# Re-run of the LR-pair table for Fgf1 -> Fgfr1: the "Sending cluster"/
# "Target cluster" set now includes ECs (previously Target cluster only
# ran over FAPs/sCs/ECs while Sending cluster only ran over FAPs/sCs -
# ECs is now also used as a sending cluster), giving the full 3x3 grid
# of (ECs, FAPs, sCs) x (ECs, FAPs, sCs) = 9 data rows instead of 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf1"
$ws.Range("C2").Value = "Fgfr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 2.004760666666666
$ws.Range("H2").Value = 6.014282
$ws.Range("I2").Value = 0.1200698528618338
$ws.Range("J2").Value = 0.1200698528618338
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.675378666666666
$ws.Range("N2").Value = 14.026136
$ws.Range("O2").Value = 0.03681964474327726
$ws.Range("P2").Value = 0.03681964474327726
$ws.Range("Q2").Value = 9.373015252705775
$ws.Range("R2").Value = 84.35713727435198
$ws.Range("S2").Value = 0.004420929326750292
$ws.Range("T2").Value = 0.004420929326750292

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf1"
$ws.Range("C3").Value = "Fgfr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 2.004760666666666
$ws.Range("H3").Value = 6.014282
$ws.Range("I3").Value = 0.1200698528618338
$ws.Range("J3").Value = 0.1200698528618338
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 82.95722966666666
$ws.Range("N3").Value = 248.871689
$ws.Range("O3").Value = 0.653306596744776
$ws.Range("P3").Value = 0.653306596744776
$ws.Range("Q3").Value = 166.3093910513664
$ws.Range("R3").Value = 1496.784519462298
$ws.Range("S3").Value = 0.07844242694481063
$ws.Range("T3").Value = 0.07844242694481063

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf1"
$ws.Range("C4").Value = "Fgfr1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 2.004760666666666
$ws.Range("H4").Value = 6.014282
$ws.Range("I4").Value = 0.1200698528618338
$ws.Range("J4").Value = 0.1200698528618338
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 39.34793966666667
$ws.Range("N4").Value = 118.043819
$ws.Range("O4").Value = 0.3098737585119468
$ws.Range("P4").Value = 0.3098737585119468
$ws.Range("Q4").Value = 78.88320175810644
$ws.Range("R4").Value = 709.948815822958
$ws.Range("S4").Value = 0.03720649659027286
$ws.Range("T4").Value = 0.03720649659027286

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf1"
$ws.Range("C5").Value = "Fgfr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.409654999999999
$ws.Range("H5").Value = 28.228965
$ws.Range("I5").Value = 0.5635664696121425
$ws.Range("J5").Value = 0.5635664696121424
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.675378666666666
$ws.Range("N5").Value = 14.026136
$ws.Range("O5").Value = 0.03681964474327726
$ws.Range("P5").Value = 0.03681964474327726
$ws.Range("Q5").Value = 43.99370024769333
$ws.Range("R5").Value = 395.94330222924
$ws.Range("S5").Value = 0.02075031720034205
$ws.Range("T5").Value = 0.02075031720034205

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf1"
$ws.Range("C6").Value = "Fgfr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 9.409654999999999
$ws.Range("H6").Value = 28.228965
$ws.Range("I6").Value = 0.5635664696121425
$ws.Range("J6").Value = 0.5635664696121424
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 82.95722966666666
$ws.Range("N6").Value = 248.871689
$ws.Range("O6").Value = 0.653306596744776
$ws.Range("P6").Value = 0.653306596744776
$ws.Range("Q6").Value = 780.5989109190982
$ws.Range("R6").Value = 7025.390198271884
$ws.Range("S6").Value = 0.3681816923017771
$ws.Range("T6").Value = 0.3681816923017769

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf1"
$ws.Range("C7").Value = "Fgfr1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 9.409654999999999
$ws.Range("H7").Value = 28.228965
$ws.Range("I7").Value = 0.5635664696121425
$ws.Range("J7").Value = 0.5635664696121424
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 39.34793966666667
$ws.Range("N7").Value = 118.043819
$ws.Range("O7").Value = 0.3098737585119468
$ws.Range("P7").Value = 0.3098737585119468
$ws.Range("Q7").Value = 370.2505372241483
$ws.Range("R7").Value = 3332.254835017335
$ws.Range("S7").Value = 0.1746344601100235
$ws.Range("T7").Value = 0.1746344601100234

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf1"
$ws.Range("C8").Value = "Fgfr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.282203999999999
$ws.Range("H8").Value = 15.846612
$ws.Range("I8").Value = 0.3163636775260238
$ws.Range("J8").Value = 0.3163636775260238
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.675378666666666
$ws.Range("N8").Value = 14.026136
$ws.Range("O8").Value = 0.03681964474327726
$ws.Range("P8").Value = 0.03681964474327726
$ws.Range("Q8").Value = 24.69630389458133
$ws.Range("R8").Value = 222.266735051232
$ws.Range("S8").Value = 0.01164839821618492
$ws.Range("T8").Value = 0.01164839821618492

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf1"
$ws.Range("C9").Value = "Fgfr1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.282203999999999
$ws.Range("H9").Value = 15.846612
$ws.Range("I9").Value = 0.3163636775260238
$ws.Range("J9").Value = 0.3163636775260238
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 82.95722966666666
$ws.Range("N9").Value = 248.871689
$ws.Range("O9").Value = 0.653306596744776
$ws.Range("P9").Value = 0.653306596744776
$ws.Range("Q9").Value = 438.1970103741853
$ws.Range("R9").Value = 3943.773093367668
$ws.Range("S9").Value = 0.2066824774981884
$ws.Range("T9").Value = 0.2066824774981884

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf1"
$ws.Range("C10").Value = "Fgfr1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.282203999999999
$ws.Range("H10").Value = 15.846612
$ws.Range("I10").Value = 0.3163636775260238
$ws.Range("J10").Value = 0.3163636775260238
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 39.34793966666667
$ws.Range("N10").Value = 118.043819
$ws.Range("O10").Value = 0.3098737585119468
$ws.Range("P10").Value = 0.3098737585119468
$ws.Range("Q10").Value = 207.8438442990253
$ws.Range("R10").Value = 1870.594598691228
$ws.Range("S10").Value = 0.09803280181165051
$ws.Range("T10").Value = 0.0980328018116505
